# Reading the file, validating data, writing the data.
# Replace the roster with the validated/cleaned employee dataset and
# drop the hyperlinks for the two rows whose e-mail cell is now blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Hyperlinks - rebuild the collection: the validated data set no
#    longer has an e-mail for two of the rows (old B10 / B8), so those
#    links disappear; the rest keep pointing at the same (cleaned)
#    addresses.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B16"), "mailto:Fifteen@fiv.com")   | Out-Null
$ws.Hyperlinks.Add($ws.Range("B15"), "mailto:Fourteen@one.com")  | Out-Null
$ws.Hyperlinks.Add($ws.Range("B14"), "mailto:Thirteen@cu.com")   | Out-Null
$ws.Hyperlinks.Add($ws.Range("B13"), "mailto:Abc123@arw.com")    | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:Ten@c.com")         | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"),  "mailto:emp_one@a.com")     | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"),  "mailto:N3@")                | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"),  "mailto:Five@.com")         | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"),  "mailto:Eight@.com")        | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"),  "mailto:empTwo@e.com")      | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"),  "mailto:four@.com")         | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"),  "mailto:XX@.com")           | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:Eight@.com")        | Out-Null

# ---------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Salary"
$ws.Range("E1").Value = "Department"

# ---------------------------------------------------------------------
# 3. Data rows. A handful of "Salary"/"Age" cells hold numeric-looking
#    text (validation artefacts) - those get an explicit Text format
#    *before* the value is written so they are not coerced back into
#    numbers.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Employee one"
$ws.Range("B2").Value = "emp_one@a.com"
$ws.Range("C2").Value = "twenty"
$ws.Range("D2").Value = 15000
$ws.Range("E2").Value = "Dept 1"

$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "empTwo@e.com"
$ws.Range("C3").Value = 55
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "60000"
$ws.Range("E3").Value = "Dept 2"

$ws.Range("A4").Value = "Employee three"
$ws.Range("B4").Value = "N3@"
$ws.Range("C4").Value = 35
$ws.Range("D4").Value = 40000
$ws.Range("E4").Value = "Dept 3"

$ws.Range("A5").Value = "Employee four"
$ws.Range("B5").Value = "four@.com"
$ws.Range("C5").Value = 20
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "50000.00"
$ws.Range("E5").ClearContents()

$ws.Range("A6").Value = "Employee five"
$ws.Range("B6").Value = "Five@.com"
$ws.Range("C6").Value = 55
$ws.Range("D6").Value = 60000
$ws.Range("E6").Value = "Dept 3"

$ws.Range("A7").Value = "Employee six"
$ws.Range("B7").Value = "XX@.com"
$ws.Range("C7").Value = 76
$ws.Range("D7").Value = 60000
$ws.Range("E7").Value = "Dept 3"

$ws.Range("A8").Value = "Employee seven"
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = 33
$ws.Range("D8").Value = 60000
$ws.Range("E8").ClearContents()

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "8"
$ws.Range("B9").Value = "Eight@.com"
$ws.Range("C9").Value = 60
$ws.Range("D9").Value = 15000
$ws.Range("E9").Value = "Dept 2"

$ws.Range("A10").Value = "Employee nine"
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = 45
$ws.Range("D10").Value = 23456
$ws.Range("E10").Value = "Dept 3"

$ws.Range("A11").Value = "Employee ten"
$ws.Range("B11").Value = "Ten@c.com"
$ws.Range("C11").Value = 65
$ws.Range("D11").Value = 12443
$ws.Range("E11").Value = "Dept 4"

$ws.Range("A12").Value = "Employee eleven"
$ws.Range("B12").Value = "Eight@.com"
$ws.Range("C12").Value = 44
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "80.00"
$ws.Range("E12").Value = "Dept 4"

$ws.Range("A13").Value = "Employee twelve"
$ws.Range("B13").Value = "Abc123@arw.com"
$ws.Range("C13").Value = 32
$ws.Range("D13").Value = 60000
$ws.Range("E13").Value = "Dept 3"

$ws.Range("A14").Value = "Employee 13"
$ws.Range("B14").Value = "Thirteen@cu.com"
$ws.Range("C14").Value = 33
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60000"
$ws.Range("E14").Value = "depart4"

$ws.Range("A15").Value = "Employee fourteen"
$ws.Range("B15").Value = "Fourteen@one.com"
$ws.Range("C15").Value = 22
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5248"
$ws.Range("E15").Value = "Dept 1"

$ws.Range("A16").Value = "Employee fiften"
$ws.Range("B16").Value = "Fifteen@fiv.com"
$ws.Range("C16").Value = 45
$ws.Range("D16").Value = 60000
$ws.Range("E16").Value = "Dept 1"

# ---------------------------------------------------------------------
# 4. Apply "Text" number format + center alignment to the whole table
#    now that every value is in place - this only touches *display*
#    formatting, the numeric cells above stay numbers.
# ---------------------------------------------------------------------
$ws.Range("A1:E1").NumberFormat = "@"

$dataFmt = $ws.Range("A2:E16")
$dataFmt.NumberFormat = "@"
$dataFmt.HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------
# 5. Column widths / sheet view - columns A and D grew a little to fit
#    the new longest entries ("Employee eleven" / "50000.00").
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 12.43
$ws.Columns("D").ColumnWidth = 9.43

$ws.Range("E12").Select()
